$d = $word.ActiveDocument

# --- Edit 1: title paragraph "RTEMS su RPi per applicazioni  “aerospace-like”"
#     becomes two centered paragraphs:
#       "RTEMS su Raspberry Pi"
#       " per applicazioni real-time"
$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2A93DB79" w14:textId="0BC54D44" w:rsidR="00C24072" w:rsidRPr="00075D84" w:rsidRDefault="003071AE" w:rsidP="00681607"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="52"/><w:szCs w:val="52"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00075D84"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="52"/><w:szCs w:val="52"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t xml:space="preserve">RTEMS su </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00075D84"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="52"/><w:szCs w:val="52"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>R</w:t></w:r><w:r w:rsidRPr="00075D84"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="52"/><w:szCs w:val="52"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>aspberry</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00075D84"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="52"/><w:szCs w:val="52"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t xml:space="preserve"> Pi</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="52"/><w:szCs w:val="52"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00075D84"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="52"/><w:szCs w:val="52"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t xml:space="preserve"> per applicazioni real-time</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("RPi per applicazioni")
if ($found1) {
    $titlePara = $rng1.Paragraphs(1).Range
    $titlePara.InsertXML($titleXml)
}

# --- Edit 2: "Anno Accademico 2020-2021" becomes "Anno Accademico 19-2020"
#     (this is the last paragraph in the document body, so we target the
#     range up to, but excluding, its trailing paragraph mark to avoid
#     leaving a stray empty paragraph behind)
$annoXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6F08E1F5" w14:textId="2CCD763B" w:rsidR="00831DDF" w:rsidRPr="00075D84" w:rsidRDefault="002F5B48" w:rsidP="00F7100E"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="00075D84"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Anno A</w:t></w:r><w:r w:rsidR="00C24072" w:rsidRPr="00075D84"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">ccademico </w:t></w:r><w:r w:rsidR="003071AE" w:rsidRPr="00075D84"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>19-2020</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Anno Accademico")
if ($found2) {
    $annoPara = $rng2.Paragraphs(1).Range
    $innerRange = $d.Range($annoPara.Start, $annoPara.End - 1)
    $innerRange.InsertXML($annoXml)
}
